$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the series. Insert a row at 74;
# Excel shifts the existing rows 74:322 down to 75:323 and grows the used
# range to A1:R323 for us.
$ws.Rows("74").Insert()

# Fill in the newly inserted row with the new record. All of the
# non-varying descriptive columns (market/region/category/etc.) share the
# same values as every other row in this sheet.
$ws.Range("A74").Value = 3
$ws.Range("B74").Value = "Femacal de La Calera"
$ws.Range("C74").Value = "Coquimbo"
$ws.Range("D74").Value = [DateTime]"2022-06-03"
$ws.Range("E74").Value = 5
$ws.Range("F74").Value = 100112039
$ws.Range("G74").Value = "Ciboulette"
$ws.Range("H74").Value = "Sin especificar"
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 138
$ws.Range("K74").Value = 1500
$ws.Range("L74").Value = 1500
$ws.Range("M74").Value = 1500
$ws.Range("N74").Value = "$/docena de atados"
$ws.Range("O74").Value = "Provincia de Quillota"
$ws.Range("P74").Value = 500
$ws.Range("Q74").Value = 3
$ws.Range("R74").Value = "Hortaliza"
